$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")
$ws3 = $wb.Worksheets.Item("Summary")

# Ensure the currency-text cell on the Summary sheet keeps its text type
# instead of being auto-converted to a number by Excel.
$ws3.Cells.Item(2, 2).NumberFormat = "@"

# --- Sheet: Top 50 Cryptocurrencies ---
$ws1.Cells.Item(2, 3).Value = 98908
$ws1.Cells.Item(2, 4).Value = 1959097766308
$ws1.Cells.Item(2, 5).Value = 117523219357
$ws1.Cells.Item(2, 6).Value = 1.99672
$ws1.Cells.Item(3, 3).Value = 3368.9
$ws1.Cells.Item(3, 4).Value = 405800320131
$ws1.Cells.Item(3, 5).Value = 57542450103
$ws1.Cells.Item(3, 6).Value = 8.200049999999999
$ws1.Cells.Item(4, 3).Value = 1
$ws1.Cells.Item(4, 4).Value = 130852749175
$ws1.Cells.Item(4, 5).Value = 130067665585
$ws1.Cells.Item(4, 6).Value = -0.22882
$ws1.Cells.Item(5, 3).Value = 262.46
$ws1.Cells.Item(5, 4).Value = 124608104880
$ws1.Cells.Item(5, 5).Value = 15419995084
$ws1.Cells.Item(5, 6).Value = 10.2062
$ws1.Cells.Item(6, 3).Value = 635.28
$ws1.Cells.Item(6, 4).Value = 92807675426
$ws1.Cells.Item(6, 5).Value = 2516065309
$ws1.Cells.Item(6, 6).Value = 4.67069
$ws1.Cells.Item(7, 3).Value = 1.41
$ws1.Cells.Item(7, 4).Value = 79540572467
$ws1.Cells.Item(7, 5).Value = 17360543557
$ws1.Cells.Item(7, 6).Value = 27.87191
$ws1.Cells.Item(8, 3).Value = 0.39439
$ws1.Cells.Item(8, 4).Value = 58022297417
$ws1.Cells.Item(8, 5).Value = 10350503719
$ws1.Cells.Item(8, 6).Value = 3.64119
$ws1.Cells.Item(9, 3).Value = 0.999438
$ws1.Cells.Item(9, 4).Value = 38245803318
$ws1.Cells.Item(9, 5).Value = 15337728230
$ws1.Cells.Item(9, 6).Value = -0.36921
$ws1.Cells.Item(10, 3).Value = 3369.5
$ws1.Cells.Item(10, 4).Value = 33013322092
$ws1.Cells.Item(10, 5).Value = 148862892
$ws1.Cells.Item(10, 6).Value = 8.480880000000001
$ws1.Cells.Item(11, 3).Value = 0.875708
$ws1.Cells.Item(11, 4).Value = 31240593922
$ws1.Cells.Item(11, 5).Value = 3600505014
$ws1.Cells.Item(11, 6).Value = 11.54925
$ws1.Cells.Item(12, 3).Value = 0.200701
$ws1.Cells.Item(12, 4).Value = 17345913079
$ws1.Cells.Item(12, 5).Value = 1083382796
$ws1.Cells.Item(12, 6).Value = 2.01243
$ws1.Cells.Item(13, 3).Value = 36.31
$ws1.Cells.Item(13, 4).Value = 14847880210
$ws1.Cells.Item(13, 5).Value = 1051213844
$ws1.Cells.Item(13, 6).Value = 8.219150000000001
$ws1.Cells.Item(14, 3).Value = 0.00002497
$ws1.Cells.Item(14, 4).Value = 14731569383
$ws1.Cells.Item(14, 5).Value = 1612892764
$ws1.Cells.Item(14, 6).Value = 4.98286
$ws1.Cells.Item(15, 3).Value = 98573
$ws1.Cells.Item(15, 4).Value = 14416601502
$ws1.Cells.Item(15, 5).Value = 906539844
$ws1.Cells.Item(15, 6).Value = 2.03488
$ws1.Cells.Item(16, 3).Value = 3986.68
$ws1.Cells.Item(16, 4).Value = 14404222176
$ws1.Cells.Item(16, 5).Value = 170561570
$ws1.Cells.Item(16, 6).Value = 8.663069999999999
$ws1.Cells.Item(17, 4).Value = 14132776098
$ws1.Cells.Item(17, 5).Value = 620294970
$ws1.Cells.Item(17, 6).Value = 4.18105
$ws1.Cells.Item(18, 3).Value = 3.63
$ws1.Cells.Item(18, 4).Value = 10344756784
$ws1.Cells.Item(18, 5).Value = 2290801051
$ws1.Cells.Item(18, 6).Value = 3.47325
$ws1.Cells.Item(19, 3).Value = 492.32
$ws1.Cells.Item(19, 4).Value = 9747842327
$ws1.Cells.Item(19, 5).Value = 2304543250
$ws1.Cells.Item(19, 6).Value = 5.4817
$ws1.Cells.Item(20, 3).Value = 3369.49
$ws1.Cells.Item(20, 4).Value = 9689653263
$ws1.Cells.Item(20, 5).Value = 2280080043
$ws1.Cells.Item(20, 6).Value = 8.377230000000001
$ws1.Cells.Item(21, 3).Value = 15.1
$ws1.Cells.Item(21, 4).Value = 9464894293
$ws1.Cells.Item(21, 5).Value = 1224865617
$ws1.Cells.Item(21, 6).Value = 5.4432
$ws1.Cells.Item(22, 3).Value = 0.00002143
$ws1.Cells.Item(22, 4).Value = 9017977603
$ws1.Cells.Item(22, 5).Value = 7053073138
$ws1.Cells.Item(22, 6).Value = 13.76133
$ws1.Cells.Item(23, 3).Value = 6.18
$ws1.Cells.Item(23, 4).Value = 8890725340
$ws1.Cells.Item(23, 5).Value = 813834898
$ws1.Cells.Item(23, 6).Value = 9.62811
$ws1.Cells.Item(24, 3).Value = 0.287069
$ws1.Cells.Item(24, 4).Value = 8573585584
$ws1.Cells.Item(24, 5).Value = 2350798434
$ws1.Cells.Item(24, 6).Value = 20.15183
$ws1.Cells.Item(25, 3).Value = 8.74
$ws1.Cells.Item(25, 4).Value = 8086736861
$ws1.Cells.Item(25, 5).Value = 3486171
$ws1.Cells.Item(25, 6).Value = 2.21077
$ws1.Cells.Item(26, 4).Value = 7044182679
$ws1.Cells.Item(26, 5).Value = 1020700272
$ws1.Cells.Item(26, 6).Value = 5.52506
$ws1.Cells.Item(27, 3).Value = 90.17
$ws1.Cells.Item(27, 4).Value = 6785098387
$ws1.Cells.Item(27, 5).Value = 1470715698
$ws1.Cells.Item(27, 6).Value = 5.87394
$ws1.Cells.Item(28, 4).Value = 6453178954
$ws1.Cells.Item(28, 5).Value = 896912955
$ws1.Cells.Item(28, 6).Value = 5.1022
$ws1.Cells.Item(29, 3).Value = 3550.07
$ws1.Cells.Item(29, 4).Value = 6115781119
$ws1.Cells.Item(29, 5).Value = 99952442
$ws1.Cells.Item(29, 6).Value = 8.49728
$ws1.Cells.Item(30, 4).Value = 5618106680
$ws1.Cells.Item(30, 5).Value = 843765851
$ws1.Cells.Item(30, 6).Value = 7.27596
$ws1.Cells.Item(31, 3).Value = 0.195965
$ws1.Cells.Item(31, 4).Value = 5329889908
$ws1.Cells.Item(31, 5).Value = 115257302
$ws1.Cells.Item(31, 6).Value = 10.52595
$ws1.Cells.Item(32, 3).Value = 0.999322
$ws1.Cells.Item(32, 4).Value = 5247228065
$ws1.Cells.Item(32, 5).Value = 15845976
$ws1.Cells.Item(32, 6).Value = -0.27986
$ws1.Cells.Item(33, 3).Value = 0.129255
$ws1.Cells.Item(33, 4).Value = 4971045961
$ws1.Cells.Item(33, 5).Value = 837343150
$ws1.Cells.Item(33, 6).Value = 2.23502
$ws1.Cells.Item(34, 3).Value = 9.640000000000001
$ws1.Cells.Item(34, 4).Value = 4574014220
$ws1.Cells.Item(34, 5).Value = 271213610
$ws1.Cells.Item(34, 6).Value = 8.047980000000001
$ws1.Cells.Item(35, 3).Value = 27.84
$ws1.Cells.Item(35, 4).Value = 4162099848
$ws1.Cells.Item(35, 5).Value = 908771366
$ws1.Cells.Item(35, 6).Value = 6.91992
$ws1.Cells.Item(36, 3).Value = 0.00005231
$ws1.Cells.Item(36, 4).Value = 3935128733
$ws1.Cells.Item(36, 5).Value = 1815144801
$ws1.Cells.Item(36, 6).Value = 7.91315
$ws1.Cells.Item(37, 3).Value = 7.41
$ws1.Cells.Item(37, 4).Value = 3836472884
$ws1.Cells.Item(37, 5).Value = 443348294
$ws1.Cells.Item(37, 6).Value = 1.65585
$ws1.Cells.Item(38, 3).Value = 0.150396
$ws1.Cells.Item(38, 4).Value = 3792375092
$ws1.Cells.Item(38, 5).Value = 155813335
$ws1.Cells.Item(38, 6).Value = 0.70825
$ws1.Cells.Item(39, 3).Value = 509.73
$ws1.Cells.Item(39, 4).Value = 3763333252
$ws1.Cells.Item(39, 5).Value = 288200562
$ws1.Cells.Item(39, 6).Value = 4.82107
$ws1.Cells.Item(40, 3).Value = 0.467259
$ws1.Cells.Item(40, 4).Value = 3724428423
$ws1.Cells.Item(40, 5).Value = 473813459
$ws1.Cells.Item(40, 6).Value = 8.26056
$ws1.Cells.Item(41, 4).Value = 3687567717
$ws1.Cells.Item(41, 5).Value = 237089640
$ws1.Cells.Item(41, 6).Value = -0.17128
$ws1.Cells.Item(42, 3).Value = 24.78
$ws1.Cells.Item(42, 4).Value = 3573423906
$ws1.Cells.Item(42, 5).Value = 38290921
$ws1.Cells.Item(42, 6).Value = 2.73603
$ws1.Cells.Item(43, 3).Value = 3.91
$ws1.Cells.Item(43, 4).Value = 3529178084
$ws1.Cells.Item(43, 5).Value = 306842077
$ws1.Cells.Item(43, 6).Value = 6.83721
$ws1.Cells.Item(44, 3).Value = 0.999465
$ws1.Cells.Item(44, 4).Value = 3441616760
$ws1.Cells.Item(44, 5).Value = 163474938
$ws1.Cells.Item(44, 6).Value = -0.28609
$ws1.Cells.Item(45, 1).Value = 'dogwifhat'
$ws1.Cells.Item(45, 2).Value = 'wif'
$ws1.Cells.Item(45, 3).Value = 3.37
$ws1.Cells.Item(45, 4).Value = 3367531542
$ws1.Cells.Item(45, 5).Value = 1269991304
$ws1.Cells.Item(45, 6).Value = 8.123480000000001
$ws1.Cells.Item(46, 1).Value = 'Artificial Superintelligence Alliance'
$ws1.Cells.Item(46, 2).Value = 'fet'
$ws1.Cells.Item(46, 3).Value = 1.28
$ws1.Cells.Item(46, 4).Value = 3343687076
$ws1.Cells.Item(46, 5).Value = 496379708
$ws1.Cells.Item(46, 6).Value = 4.92834
$ws1.Cells.Item(47, 3).Value = 0.773817
$ws1.Cells.Item(47, 4).Value = 3171101048
$ws1.Cells.Item(47, 5).Value = 1667534686
$ws1.Cells.Item(47, 6).Value = 14.03529
$ws1.Cells.Item(48, 3).Value = 160.25
$ws1.Cells.Item(48, 4).Value = 2957759868
$ws1.Cells.Item(48, 5).Value = 84034592
$ws1.Cells.Item(48, 6).Value = -1.14394
$ws1.Cells.Item(49, 3).Value = 1.94
$ws1.Cells.Item(49, 4).Value = 2918415642
$ws1.Cells.Item(49, 5).Value = 422364809
$ws1.Cells.Item(49, 6).Value = 2.84764
$ws1.Cells.Item(50, 1).Value = 'Filecoin'
$ws1.Cells.Item(50, 2).Value = 'fil'
$ws1.Cells.Item(50, 3).Value = 4.67
$ws1.Cells.Item(50, 4).Value = 2805339303
$ws1.Cells.Item(50, 5).Value = 584722236
$ws1.Cells.Item(50, 6).Value = 8.904389999999999
$ws1.Cells.Item(51, 1).Value = 'OKB'
$ws1.Cells.Item(51, 2).Value = 'okb'
$ws1.Cells.Item(51, 3).Value = 46.56
$ws1.Cells.Item(51, 4).Value = 2795155174
$ws1.Cells.Item(51, 5).Value = 20372276
$ws1.Cells.Item(51, 6).Value = 6.04536

# --- Sheet: Top 5 by Market Cap ---
$ws2.Cells.Item(2, 2).Value = 1959097766308
$ws2.Cells.Item(3, 2).Value = 405800320131
$ws2.Cells.Item(4, 2).Value = 130852749175
$ws2.Cells.Item(5, 2).Value = 124608104880
$ws2.Cells.Item(6, 2).Value = 92807675426

# --- Sheet: Summary ---
$ws3.Cells.Item(2, 2).Value = '$4350.40'
$ws3.Cells.Item(3, 2).Value = 'XRP (27.87%)'
$ws3.Cells.Item(4, 2).Value = 'Monero (-1.14%)'
Write-Host "Applied crypto live data update."
